$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZWNL-Treinseries 2021")
$rng = $ws.Range("A2:F29")
$rng.Sort($ws.Range("A2:A29"))
